$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated duty-roster assignments for column B (rows 2-31).
$values = @{
    2  = "豊島亮"
    3  = "兒島大志郎"
    4  = "山口玲, Cox Matthew Jonah"
    5  = "日高泰聖"
    6  = "志塚惇希"
    7  = "白岩詩佑介"
    8  = "富澤天音"
    9  = "石井海成"
    10 = "Nicholas Tristan Aryasatyo"
    11 = "小溝賢"
    12 = "小野文哉"
    13 = "渡部魁"
    14 = "崎谷航平"
    16 = "三神佳誠"
    17 = "氏家琉貴, Hansen Jakob U"
    18 = "羽賀尚生"
    19 = "島田実"
    20 = "足立耕平"
    21 = "遠藤隼人"
    22 = "富澤天音"
    23 = "神山修造"
    24 = "川田涼介"
    25 = "豊島亮"
    26 = "兒島大志郎"
    27 = "山口玲"
    28 = "日高泰聖"
    29 = "氏家琉貴"
    30 = "志塚惇希"
    31 = "白岩詩佑介"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}

# Selected cell moves to G21 in the saved view.
$ws.Range("G21").Select()
